# Remarque - Problematique.xlsx
# Commit: "afficher image selon User src" — add a new tracker row (row 35)
# documenting the change-password controller / getDoctrine() issue and its
# StackOverflow fix link, mirroring the formatting of the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use row 26 as the formatting template (same style pattern the new row needs:
# plain body style for A/C/D/F/G, filled style for B, date style for E, and the
# hyperlink style for H) and copy it down into the new row 35.
$ws.Range("A26:H26").Copy()
$ws.Range("A35:H35").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Row height for the new entry
$ws.Rows.Item(35).RowHeight = 120

# Cell values
$ws.Range("B35").Value = "change-password controller ne connait pas ""em"""
$ws.Range("D35").Value = "public function getDoctrine()" + [char]10 + "{" + [char]10 + "    return `$this->container->get('doctrine');" + [char]10 + "}"
$ws.Range("E35").Value = 42107
$ws.Range("H35").Value = "http://stackoverflow.com/questions/22841440/symfony2-call-to-undefined-method-getdoctrine-when-overriding-fosuserbundle"

# make sure the cells that should stay empty really are
$ws.Range("A35").Value = ""
$ws.Range("C35").Value = ""
$ws.Range("F35").Value = ""
$ws.Range("G35").Value = ""

# Hyperlink on H35 pointing to the StackOverflow answer
$ws.Hyperlinks.Add($ws.Range("H35"), "http://stackoverflow.com/questions/22841440/symfony2-call-to-undefined-method-getdoctrine-when-overriding-fosuserbundle")

# Keep the sheet's dimension/view in sync with the newly-used row, same as
# Excel does automatically when the used range grows.
$ws.Application.ActiveWindow.ScrollRow = 34
$ws.Range("B35").Select()
